# Insert a new weekly price record for "Espinaca" (Vega Central Mapocho de
# Santiago) as row 335, pushing the former rows 335-394 down to 336-395.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 335..394 down by one row.
$ws.Rows("335:335").Insert()

# Populate the newly inserted row 335 with the new record.
$ws.Cells.Item(335, 1).Value  = 9
$ws.Cells.Item(335, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(335, 3).Value  = "Metropolitana"
$ws.Cells.Item(335, 4).Value  = 44641
$ws.Cells.Item(335, 5).Value  = 13
$ws.Cells.Item(335, 6).Value  = 100112012
$ws.Cells.Item(335, 7).Value  = "Espinaca"
$ws.Cells.Item(335, 8).Value  = "Sin especificar"
$ws.Cells.Item(335, 9).Value  = "Primera"
$ws.Cells.Item(335, 10).Value = 61
$ws.Cells.Item(335, 11).Value = 12000
$ws.Cells.Item(335, 12).Value = 14000
$ws.Cells.Item(335, 13).Value = 13016
$ws.Cells.Item(335, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(335, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(335, 16).Value = 1302
$ws.Cells.Item(335, 17).Value = 10
$ws.Cells.Item(335, 18).Value = "Hortaliza"
